# Insert a new weekly price-report row for "Piña" (Vega Modelo de Temuco)
# at sheet row 398, pushing the existing rows 398-466 down to 399-467.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 398..466 down one position (mirrors the shared layout of every
# other row for this product/market), making room for the new record.
$ws.Rows(398).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(398, 1).Value  = 10
$ws.Cells.Item(398, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(398, 3).Value  = "La Araucanía"
$ws.Cells.Item(398, 4).Value  = 44694
$ws.Cells.Item(398, 5).Value  = 9
$ws.Cells.Item(398, 6).Value  = "Fruta"
$ws.Cells.Item(398, 7).Value  = 100108
$ws.Cells.Item(398, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(398, 9).Value  = 100108005
$ws.Cells.Item(398, 10).Value = "Piña"
$ws.Cells.Item(398, 11).Value = "Caramelo"
$ws.Cells.Item(398, 12).Value = "Primera"
$ws.Cells.Item(398, 13).Value = 55
$ws.Cells.Item(398, 14).Value = 19000
$ws.Cells.Item(398, 15).Value = 19000
$ws.Cells.Item(398, 16).Value = 19000
$ws.Cells.Item(398, 17).Value = "`$/caja 12 unidades"
$ws.Cells.Item(398, 18).Value = "Ecuador"
$ws.Cells.Item(398, 19).Value = 1583
$ws.Cells.Item(398, 20).Value = 12
